$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1851851851851852
$ws.Range("C2").Value = 0.5343915343915344
$ws.Range("J2").Value = 0.01058201058201058
$ws.Range("P2").Value = 0.1693121693121693
$ws.Range("S2").Value = 0.1005291005291005
$ws.Range("C3").Value = 0.00980392156862745
$ws.Range("J3").Value = 0.04901960784313725
$ws.Range("P3").Value = 0.7450980392156863
$ws.Range("S3").Value = 0.196078431372549
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7073170731707317
$ws.Range("S4").Value = 0.2682926829268293
$ws.Range("B6").Value = 0.05208333333333334
$ws.Range("D6").Value = 0.015625
$ws.Range("F6").Value = 0.03645833333333334
$ws.Range("J6").Value = 0.2135416666666667
$ws.Range("O6").Value = 0.005208333333333333
$ws.Range("Q6").Value = 0.1145833333333333
$ws.Range("R6").Value = 0.0625
$ws.Range("S6").Value = 0.5
$ws.Range("B7").Value = 0.05555555555555555
$ws.Range("D7").Value = 0.04320987654320987
$ws.Range("E7").Value = 0.006172839506172839
$ws.Range("F7").Value = 0.06172839506172839
$ws.Range("J7").Value = 0.1234567901234568
$ws.Range("O7").Value = 0.01851851851851852
$ws.Range("Q7").Value = 0.1728395061728395
$ws.Range("R7").Value = 0.1358024691358025
$ws.Range("S7").Value = 0.382716049382716
$ws.Range("B8").Value = 0.06004140786749482
$ws.Range("D8").Value = 0.01449275362318841
$ws.Range("F8").Value = 0.07867494824016563
$ws.Range("J8").Value = 0.08695652173913043
$ws.Range("O8").Value = 0.02691511387163561
$ws.Range("Q8").Value = 0.1697722567287785
$ws.Range("R8").Value = 0.09109730848861283
$ws.Range("S8").Value = 0.4720496894409938
$ws.Range("B9").Value = 0.05978260869565218
$ws.Range("D9").Value = 0.005434782608695652
$ws.Range("F9").Value = 0.05978260869565218
$ws.Range("J9").Value = 0.07065217391304347
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0.01630434782608696
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0.1467391304347826
$ws.Range("R9").Value = 0.08695652173913043
$ws.Range("S9").Value = 0.5543478260869565
$ws.Range("B10").Value = 0.07949412827461608
$ws.Range("D10").Value = 0.02258355916892502
$ws.Range("F10").Value = 0.06142728093947606
$ws.Range("J10").Value = 0.1273712737127371
$ws.Range("O10").Value = 0.01084010840108401
$ws.Range("Q10").Value = 0.1996386630532972
$ws.Range("R10").Value = 0.07949412827461608
$ws.Range("S10").Value = 0.4191508581752484
$ws.Range("G11").Value = 0.1594827586206897
$ws.Range("J11").Value = 0.06465517241379311
$ws.Range("K11").Value = 0.1939655172413793
$ws.Range("L11").Value = 0.5689655172413793
$ws.Range("S11").Value = 0.01293103448275862
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1642857142857143
$ws.Range("K12").Value = 0.01428571428571429
$ws.Range("L12").Value = 0.02142857142857143
$ws.Range("S12").Value = 0.05
$ws.Range("G13").Value = 0.7428571428571429
$ws.Range("J13").Value = 0.2285714285714286
$ws.Range("S13").Value = 0.02857142857142857
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.02380952380952381
$ws.Range("H15").Value = 0.1726190476190476
$ws.Range("I15").Value = 0.125
$ws.Range("J15").Value = 0.3095238095238095
$ws.Range("K15").Value = 0.04166666666666666
$ws.Range("M15").Value = 0.01785714285714286
$ws.Range("N15").Value = 0.005952380952380952
$ws.Range("O15").Value = 0.02976190476190476
$ws.Range("S15").Value = 0.2738095238095238
$ws.Range("F16").Value = 0.007575757575757576
$ws.Range("H16").Value = 0.2424242424242424
$ws.Range("I16").Value = 0.143939393939394
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("K16").Value = 0.1060606060606061
$ws.Range("N16").Value = 0.007575757575757576
$ws.Range("O16").Value = 0.04545454545454546
$ws.Range("S16").Value = 0.1136363636363636
$ws.Range("F17").Value = 0.02419354838709677
$ws.Range("H17").Value = 0.2258064516129032
$ws.Range("I17").Value = 0.1182795698924731
$ws.Range("J17").Value = 0.3978494623655914
$ws.Range("K17").Value = 0.08333333333333333
$ws.Range("M17").Value = 0.01612903225806452
$ws.Range("O17").Value = 0.05376344086021505
$ws.Range("S17").Value = 0.08064516129032258
$ws.Range("F18").Value = 0.0111731843575419
$ws.Range("H18").Value = 0.1899441340782123
$ws.Range("I18").Value = 0.08379888268156424
$ws.Range("J18").Value = 0.4134078212290503
$ws.Range("K18").Value = 0.08379888268156424
$ws.Range("M18").Value = 0.01675977653631285
$ws.Range("O18").Value = 0.05027932960893855
$ws.Range("S18").Value = 0.1508379888268156
$ws.Range("F19").Value = 0.01105845181674566
$ws.Range("H19").Value = 0.2448657187993681
$ws.Range("I19").Value = 0.06635071090047394
$ws.Range("J19").Value = 0.391785150078989
$ws.Range("K19").Value = 0.0924170616113744
$ws.Range("M19").Value = 0.01816745655608215
$ws.Range("O19").Value = 0.05924170616113744
$ws.Range("S19").Value = 0.1161137440758294
